# benchmark now measures the time it takes each file to transmit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header label: "Time (sec)" -> "Avg Time (sec)"
$ws.Range("C1").Value = "Avg Time (sec)"

# File labels get a "(10 times)" suffix (each string is shared across 3 rows)
$ws.Range("B3").Value = "small.txt (45 Bytes) (10 times)"
$ws.Range("B6").Value = "small.txt (45 Bytes) (10 times)"
$ws.Range("B9").Value = "small.txt (45 Bytes) (10 times)"

$ws.Range("B4").Value = "medium.jpg (44.43 KB) (10 times)"
$ws.Range("B7").Value = "medium.jpg (44.43 KB) (10 times)"
$ws.Range("B10").Value = "medium.jpg (44.43 KB) (10 times)"

# Widen columns B-E (COM ColumnWidth is offset from the stored OOXML width
# by 11/14 of a character on this sheet's font, so subtract that offset to
# land exactly on the target stored widths of 30, 10, 10, 10).
$offset = 11/14
$ws.Columns.Item(2).ColumnWidth = 30 - $offset
$ws.Columns.Item(3).ColumnWidth = 10 - $offset
$ws.Columns.Item(4).ColumnWidth = 10 - $offset
$ws.Columns.Item(5).ColumnWidth = 10 - $offset

# Updated measured values (now much faster — averages of 10 runs instead of single-run timings)
$ws.Range("C3").Value = 0.0030138000147417188
$ws.Range("D3").Value = 0.002043300005607307

$ws.Range("C4").Value = 0.04863890000851825
$ws.Range("D4").Value = 0.0057667999877594415

$ws.Range("C6").Value = 0.004543900047428906
$ws.Range("D6").Value = 0.003966900007799268

$ws.Range("C7").Value = 0.07063169997418299
$ws.Range("D7").Value = 0.010418599983677269

$ws.Range("C9").Value = 0.23238089999649675
$ws.Range("D9").Value = 0.77170329997316

$ws.Range("C10").Value = 0.3117841000203043
$ws.Range("D10").Value = 0.9604642999940551
